# Lgi4-Adam11.xlsx — "update scripts wuth new tpm"
#
# The NATMI export script was re-run with refreshed TPM values. For the
# Lgi4/Adam11 ligand-receptor pair, the "ECs" sending-cluster block (the
# 3 rows where Sending cluster = ECs) was dropped entirely, and the
# remaining rows (Sending cluster = FAPs / MuSCs) were recomputed with new
# expression numbers. The sheet shrinks from 10 rows (9 data rows) down to
# 7 rows (6 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "ECs" sending-cluster block (former rows 8:10) so only the
# FAPs/MuSCs blocks remain; remaining rows shift up automatically.
$ws.Range("A8:T10").Delete() | Out-Null

# New data for rows 2-7 (columns A-T), reflecting the refreshed TPM run.
$data = @(
    @("FAPs", "Lgi4", "Adam11", "ECs",   3, 1, 10.326286,            30.978858,          0.5836270415642173, 0.5836270415642173, 2, 0.6666666666666666, 0.092903, 0.278709, 0.03600043090620505, 0.03600043090620505, 0.9593429482579999, 8.634086534322,    0.02101082498482547, 0.02101082498482547),
    @("FAPs", "Lgi4", "Adam11", "FAPs",  3, 1, 10.326286,            30.978858,          0.5836270415642173, 0.5836270415642173, 2, 0.6666666666666666, 0.292138, 0.876414, 0.1132051051535142,  0.1132051051535142,  3.016700539468,     27.150304855212,   0.06606956061071162, 0.06606956061071163),
    @("FAPs", "Lgi4", "Adam11", "MuSCs", 3, 1, 10.326286,            30.978858,          0.5836270415642173, 0.5836270415642173, 3, 1,                  2.195567, 6.586701, 0.8507944639402807,  0.8507944639402808,  22.672052774162,    204.048474967458,  0.4965466559686802,  0.4965466559686803),
    @("MuSCs","Lgi4", "Adam11", "ECs",   3, 1, 7.367010000000001,    22.10103,           0.4163729584357827, 0.4163729584357828, 2, 0.6666666666666666, 0.092903, 0.278709, 0.03600043090620505, 0.03600043090620505, 0.68441733003,      6.15975597027,     0.01498960592137958, 0.01498960592137958),
    @("MuSCs","Lgi4", "Adam11", "FAPs",  3, 1, 7.367010000000001,    22.10103,           0.4163729584357827, 0.4163729584357828, 2, 0.6666666666666666, 0.292138, 0.876414, 0.1132051051535142,  0.1132051051535142,  2.15218356738,      19.36965210642,    0.04713554454280258, 0.0471355445428026),
    @("MuSCs","Lgi4", "Adam11", "MuSCs", 3, 1, 7.367010000000001,    22.10103,           0.4163729584357827, 0.4163729584357828, 3, 1,                  2.195567, 6.586701, 0.8507944639402807,  0.8507944639402808,  16.17476404467,     145.57287640203,   0.3542478079716005,  0.3542478079716006)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}
